$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020115423449706
$ws.Range("D2").Value = 1.025312318946697
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.031246982650435
$ws.Range("I2").Value = 1.029635171504589
$ws.Range("J2").Value = 1.025314357548987
$ws.Range("K2").Value = 1.028138258276699
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.034055653289509
$ws.Range("N2").Value = 1.012487563477529

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021275462127896
$ws.Range("D3").Value = 1.026155751271619
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.032595985158203
$ws.Range("I3").Value = 1.029900930594138
$ws.Range("J3").Value = 1.026110345184985
$ws.Range("K3").Value = 1.028788868759342
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.035211730778197
$ws.Range("N3").Value = 1.012753081312279

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022025217574907
$ws.Range("D4").Value = 1.026700509672601
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.033468088916878
$ws.Range("I4").Value = 1.030070736252279
$ws.Range("J4").Value = 1.026624054819904
$ws.Range("K4").Value = 1.02920822156968
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.035958473466297
$ws.Range("N4").Value = 1.012924361282109

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022340209641418
$ws.Range("D5").Value = 1.026929287907034
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.033834535385128
$ws.Range("I5").Value = 1.030141606317767
$ws.Range("J5").Value = 1.026839697202155
$ws.Range("K5").Value = 1.029384126890153
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.036272091767464
$ws.Range("N5").Value = 1.012996241475322

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022393086238469
$ws.Range("D6").Value = 1.026967686836978
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.033896052554642
$ws.Range("I6").Value = 1.030153475454444
$ws.Range("J6").Value = 1.026875885717091
$ws.Range("K6").Value = 1.029413639290727
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.036324731429483
$ws.Range("N6").Value = 1.013008303106018

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02202942731725
$ws.Range("D7").Value = 1.026703567553438
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.033472986118942
$ws.Range("I7").Value = 1.030071685249465
$ws.Range("J7").Value = 1.026626937504063
$ws.Range("K7").Value = 1.029210573558785
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.035962665273024
$ws.Range("N7").Value = 1.012925322242934

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020507644994665
$ws.Range("D8").Value = 1.025597568184909
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.031703050231571
$ws.Range("I8").Value = 1.029725433189007
$ws.Range("J8").Value = 1.025583645213559
$ws.Range("K8").Value = 1.028358474410167
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.034446629488056
$ws.Range("N8").Value = 1.012577406038403

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.01781931706972
$ws.Range("D9").Value = 1.023640965609473
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.028577983735095
$ws.Range("I9").Value = 1.029098742072558
$ws.Range("J9").Value = 1.023734843200181
$ws.Range("K9").Value = 1.026844394350911
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.031764964706753
$ws.Range("N9").Value = 1.011960269967594

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016022398741782
$ws.Range("D10").Value = 1.022331331672458
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.026490209160944
$ws.Range("I10").Value = 1.028669788393571
$ws.Range("J10").Value = 1.022495228081145
$ws.Range("K10").Value = 1.025826488343917
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.029970132772859
$ws.Range("N10").Value = 1.011546083927289

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.015243162145317
$ws.Range("D11").Value = 1.021762989858258
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.025585087053981
$ws.Range("I11").Value = 1.028481392040306
$ws.Range("J11").Value = 1.021956759305155
$ws.Range("K11").Value = 1.025383687433847
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.029191234763901
$ws.Range("N11").Value = 1.01136607467436

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.014953541583591
$ws.Range("D12").Value = 1.021551691445549
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.025248714444976
$ws.Range("I12").Value = 1.028411013364421
$ws.Range("J12").Value = 1.021756489494399
$ws.Range("K12").Value = 1.025218903370738
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.028901654578024
$ws.Range("N12").Value = 1.011299110825484

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01501567425716
$ws.Range("D13").Value = 1.02159702433636
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.025320875273443
$ws.Range("I13").Value = 1.028426127937154
$ws.Range("J13").Value = 1.021799459780681
$ws.Range("K13").Value = 1.025254264077126
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.028963782463506
$ws.Range("N13").Value = 1.011313479351991

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.015219225657022
$ws.Range("D14").Value = 1.021745527769228
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.025557285881862
$ws.Range("I14").Value = 1.028475582678919
$ws.Range("J14").Value = 1.021940210225955
$ws.Range("K14").Value = 1.025370072625281
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.029167303358116
$ws.Range("N14").Value = 1.011360541473549

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.015344616788368
$ws.Range("D15").Value = 1.021837000289334
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.025702923606906
$ws.Range("I15").Value = 1.02850600036632
$ws.Range("J15").Value = 1.022026896901352
$ws.Range("K15").Value = 1.025441385224866
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.029292664363129
$ws.Range("N15").Value = 1.011389524676871

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016074089410323
$ws.Range("D16").Value = 1.022369023945385
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.026550255542434
$ws.Range("I16").Value = 1.028682235603331
$ws.Range("J16").Value = 1.022530928337206
$ws.Range("K16").Value = 1.025855832457847
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.030021789024856
$ws.Range("N16").Value = 1.011558016504749

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.016531355699868
$ws.Range("D17").Value = 1.02270240955679
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.027081466358819
$ws.Range("I17").Value = 1.028792071536205
$ws.Range("J17").Value = 1.022846635429226
$ws.Range("K17").Value = 1.026115256821878
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.030478685452664
$ws.Range("N17").Value = 1.011663528761811

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.016797959775698
$ws.Range("D18").Value = 1.022896746107145
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.027391206619724
$ws.Range("I18").Value = 1.028855880708849
$ws.Range("J18").Value = 1.023030617449518
$ws.Range("K18").Value = 1.026266377849019
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.030745019332062
$ws.Range("N18").Value = 1.011725008257515

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.016888846041703
$ws.Range("D19").Value = 1.022962989258798
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.027496802193599
$ws.Range("I19").Value = 1.028877594542715
$ws.Range("J19").Value = 1.02309332270836
$ws.Range("K19").Value = 1.02631787287443
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.0308358042463
$ws.Range("N19").Value = 1.011745960337352

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.016482306924101
$ws.Range("D20").Value = 1.022666653028251
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.027024483465121
$ws.Range("I20").Value = 1.028780313693116
$ws.Range("J20").Value = 1.022812780085524
$ws.Range("K20").Value = 1.026087443411702
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.030429681996781
$ws.Range("N20").Value = 1.011652214921906

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.015159289766336
$ws.Range("D21").Value = 1.02170180249235
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.025487673577252
$ws.Range("I21").Value = 1.028461030534968
$ws.Range("J21").Value = 1.021898769881144
$ws.Range("K21").Value = 1.025335978410966
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.029107378811947
$ws.Range("N21").Value = 1.01134668562841

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014326428392783
$ws.Range("D22").Value = 1.021094057133137
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.024520434678154
$ws.Range("I22").Value = 1.028257970561761
$ws.Range("J22").Value = 1.021322598578854
$ws.Range("K22").Value = 1.024861719288025
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.028274472031678
$ws.Range("N22").Value = 1.011154006042934

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.014768042497995
$ws.Range("D23").Value = 1.0214163397137
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.025033281214999
$ws.Range("I23").Value = 1.028365836091587
$ws.Range("J23").Value = 1.021628180453798
$ws.Range("K23").Value = 1.025113302549688
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.028716157241147
$ws.Range("N23").Value = 1.011256204444047

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.016504470299371
$ws.Range("D24").Value = 1.022682810240266
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.027050231909376
$ws.Range("I24").Value = 1.028785627348136
$ws.Range("J24").Value = 1.022828078366569
$ws.Range("K24").Value = 1.026100011702682
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.030451825062155
$ws.Range("N24").Value = 1.011657327356966

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01851513087441
$ws.Range("D25").Value = 1.024147711891902
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.029386647135529
$ws.Range("I25").Value = 1.029262721324723
$ws.Range("J25").Value = 1.024214043680313
$ws.Range("K25").Value = 1.027237316995462
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.012120299027044
$ws.Range("N25").Value = 1.0121202990270548
